$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F3").Value = "x-1p-user=(SYS_USER1)"
$ws.Range("F4").Value = "x-1p-user=(SYS_USER1)"
$ws.Range("F5").Value = "x-1p-user=(SYS_USER1)"
$ws.Range("F6").Value = "x-1p-user=(SYS_USER1)"
$ws.Range("F7").Value = "x-1p-user=(SYS_USER1)"
